$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58, shifting existing rows 58-124 down to 59-125.
$ws.Rows.Item(58).EntireRow.Insert()

# Populate the newly inserted row 58 with the new data record.
$ws.Cells.Item(58, 1).Value = 10
$ws.Cells.Item(58, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(58, 3).Value = "La Araucanía"
$ws.Cells.Item(58, 4).Value = 44601
$ws.Cells.Item(58, 5).Value = 9
$ws.Cells.Item(58, 6).Value = 100112012
$ws.Cells.Item(58, 7).Value = "Espinaca"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 35
$ws.Cells.Item(58, 11).Value = 12000
$ws.Cells.Item(58, 12).Value = 12000
$ws.Cells.Item(58, 13).Value = 12000
$ws.Cells.Item(58, 14).Value = "$/docena de atados"
$ws.Cells.Item(58, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(58, 16).Value = 4000
$ws.Cells.Item(58, 17).Value = 3
$ws.Cells.Item(58, 18).Value = "Hortaliza"
